$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'310.52"
$ws.Range("E2").Value = "'0.91%"
$ws.Range("G2").Value = "'12"
$ws.Range("D3").Value = "'37.36"
$ws.Range("E3").Value = "'-0.83%"
$ws.Range("G3").Value = "'12"
$ws.Range("D4").Value = "'5.103"
$ws.Range("E4").Value = "'0.22%"
$ws.Range("G4").Value = "'12"
$ws.Range("D5").Value = "'0.07787"
$ws.Range("E5").Value = "'-1.29%"
$ws.Range("G5").Value = "'12"
$ws.Range("B6").Value = "GateToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D6").Value = "'4.411"
$ws.Range("E6").Value = "'1.56%"
$ws.Range("G6").Value = "'12"
$ws.Range("B7").Value = "KuCoinToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D7").Value = "'8.205"
$ws.Range("E7").Value = "'-0.11%"
$ws.Range("G7").Value = "'12"
$ws.Range("B8").Value = "FTXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D8").Value = "'1.874"
$ws.Range("E8").Value = "'-5.55%"
$ws.Range("G8").Value = "'12"
$ws.Range("D9").Value = "'0.9279"
$ws.Range("E9").Value = "'0.42%"
$ws.Range("G9").Value = "'12"
$ws.Range("D10").Value = "'0.1175"
$ws.Range("E10").Value = "'-10.11%"
$ws.Range("G10").Value = "'12"
$ws.Range("D11").Value = "'0.1903"
$ws.Range("E11").Value = "'0.60%"
$ws.Range("G11").Value = "'12"
$ws.Range("D12").Value = "'0.09266"
$ws.Range("E12").Value = "'5.79%"
$ws.Range("G12").Value = "'12"
$ws.Range("D13").Value = "'0.03424"
$ws.Range("E13").Value = "'-0.27%"
$ws.Range("G13").Value = "'12"
$ws.Range("D14").Value = "'0.09622"
$ws.Range("E14").Value = "'-1.23%"
$ws.Range("G14").Value = "'12"
$ws.Range("D15").Value = "'0.001370"
$ws.Range("E15").Value = "'-1.90%"
$ws.Range("G15").Value = "'12"
$ws.Range("D16").Value = "'0.005891"
$ws.Range("E16").Value = "'-1.78%"
$ws.Range("G16").Value = "'12"
$ws.Range("E17").Value = "'-0.28%"
$ws.Range("G17").Value = "'12"
$ws.Range("B18").Value = "BTSEToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D18").Value = "'3.054"
$ws.Range("E18").Value = "'-1.31%"
$ws.Range("G18").Value = "'12"
$ws.Range("D19").Value = "'0.3403"
$ws.Range("E19").Value = "'-1.04%"
$ws.Range("G19").Value = "'12"
$ws.Range("D20").Value = "'5.263"
$ws.Range("E20").Value = "'5.02%"
$ws.Range("G20").Value = "'12"
$ws.Range("D21").Value = "'0.1279"
$ws.Range("E21").Value = "'-0.43%"
$ws.Range("G21").Value = "'12"
$ws.Range("D22").Value = "'0.2590"
$ws.Range("E22").Value = "'2.87%"
$ws.Range("G22").Value = "'12"
$ws.Range("E23").Value = "'180.29%"
$ws.Range("G23").Value = "'12"
$ws.Range("D24").Value = "'0.04328"
$ws.Range("E24").Value = "'-0.17%"
$ws.Range("G24").Value = "'12"
$ws.Range("D25").Value = "'0.001197"
$ws.Range("E25").Value = "'-1.84%"
$ws.Range("G25").Value = "'12"
$ws.Range("E26").Value = "'-7.75%"
$ws.Range("G26").Value = "'12"
$ws.Range("D27").Value = "'0.0001299"
$ws.Range("E27").Value = "'-63.83%"
$ws.Range("G27").Value = "'12"
$ws.Range("G28").Value = "'12"
$ws.Range("G29").Value = "'12"
$ws.Range("G30").Value = "'12"
$ws.Range("G31").Value = "'12"
$ws.Range("G32").Value = "'12"
$ws.Range("G33").Value = "'12"
$ws.Range("G34").Value = "'12"
$ws.Range("G35").Value = "'12"
$ws.Range("G36").Value = "'12"
$ws.Range("G37").Value = "'12"
$ws.Range("G38").Value = "'12"
$ws.Range("D39").Value = "'0.02069"
$ws.Range("E39").Value = "'-9.54%"
$ws.Range("G39").Value = "'12"
$ws.Range("D40").Value = "'0.05046"
$ws.Range("E40").Value = "'1.23%"
$ws.Range("G40").Value = "'12"
$ws.Range("D41").Value = "'0.007673"
$ws.Range("G41").Value = "'12"
$ws.Range("D42").Value = "'0.009785"
$ws.Range("E42").Value = "'-1.58%"
$ws.Range("G42").Value = "'12"
$ws.Range("D43").Value = "'0.1345"
$ws.Range("E43").Value = "'-0.82%"
$ws.Range("G43").Value = "'12"
$ws.Range("D44").Value = "'0.002009"
$ws.Range("E44").Value = "'-4.08%"
$ws.Range("G44").Value = "'12"
$ws.Range("D45").Value = "'0.008600"
$ws.Range("E45").Value = "'7.14%"
$ws.Range("G45").Value = "'12"
$ws.Range("D46").Value = "'0.00006708"
$ws.Range("E46").Value = "'4.78%"
$ws.Range("G46").Value = "'12"
$ws.Range("E47").Value = "'-0.21%"
$ws.Range("G47").Value = "'12"
$ws.Range("B48").Value = "CoinbaseStockToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D48").Value = "'0.001200"
$ws.Range("E48").Value = "'-0.29%"
$ws.Range("G48").Value = "'12"
$ws.Range("B49").Value = "BOLO"
$ws.Range("C49").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D49").Value = "'0.002942"
$ws.Range("E49").Value = "'-2.04%"
$ws.Range("G49").Value = "'12"
$ws.Range("D50").Value = "'0.00002100"
$ws.Range("E50").Value = "'-0.21%"
$ws.Range("G50").Value = "'12"
$ws.Range("D51").Value = "'0.0002000"
$ws.Range("E51").Value = "'-0.21%"
$ws.Range("G51").Value = "'12"
